$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.655.05"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "1.988.28"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'245.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("D7").Value = "'59.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -10.06%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").Value = "'0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.64%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "'24.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.92%  "
$ws.Range("D14").Value = "'0.873"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "'14.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").Value = "2.278.10"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").Value = "'5.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "1.985.37"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("D19").Value = "36.558.72"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "'71.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "'5.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").Value = "'235.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").Value = "'2.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").Value = "'162.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'19.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "'0.129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.90%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'1.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "'0.0631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("D40").Value = "'3.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'0.0978"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").Value = "'0.0216"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").Value = "'16.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").Value = "'93.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "1.376.13"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "'2.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "'45.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.57%  "
